# Apply updated crypto price / 1h-volume figures to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# row => @{ D = "new price text"; E = "new volume text" }
$updates = @{
    2 = @{ D='61.022.94'; E='  +1.21%  ' }
    3 = @{ D='3.384.63'; E='  +0.04%  ' }
    4 = @{ E='  +0.02%  ' }
    5 = @{ D='571.16'; E='  +0.56%  ' }
    6 = @{ D='140.97'; E='  +0.52%  ' }
    7 = @{ E='  -0.01%  ' }
    8 = @{ D='0.473'; E='  +0.33%  ' }
    9 = @{ E='  +2.76%  ' }
    10 = @{ D='0.123'; E='  -0.94%  ' }
    12 = @{ D='3.962.92'; E='  +0.05%  ' }
    14 = @{ D='27.85'; E='  -0.32%  ' }
    15 = @{ D='3.409.06'; E='  +0.43%  ' }
    16 = @{ E='  +0.33%  ' }
    17 = @{ D='61.130.33'; E='  +1.17%  ' }
    18 = @{ D='6.12'; E='  -1.21%  ' }
    19 = @{ D='13.62'; E='  -2.19%  ' }
    20 = @{ E='  -1.23%  ' }
    21 = @{ D='381.89'; E='  -1.07%  ' }
    22 = @{ D='75.90'; E='  +3.64%  ' }
    23 = @{ E='  -1.26%  ' }
    24 = @{ D='1.00'; E='  +0.02%  ' }
    25 = @{ E='  -0.32%  ' }
    26 = @{ D='3.519.68'; E='  -0.26%  ' }
    27 = @{ D='0.190'; E='  +6.77%  ' }
    28 = @{ E='  -0.05%  ' }
    29 = @{ D='7.26'; E='  -1.20%  ' }
    30 = @{ E='  +0.62%  ' }
    31 = @{ E='  +0.39%  ' }
    33 = @{ E='  -3.17%  ' }
    34 = @{ D='23.29'; E='  -1.44%  ' }
    35 = @{ D='6.96'; E='  +1.06%  ' }
    36 = @{ D='166.64'; E='  -0.58%  ' }
    37 = @{ D='3.419.63'; E='  +0.17%  ' }
    38 = @{ D='4.99'; E='  +1.74%  ' }
    39 = @{ E='  -2.55%  ' }
    40 = @{ E='  -0.25%  ' }
    41 = @{ D='26.25'; E='  -1.93%  ' }
    42 = @{ D='1.00'; E='  +0.02%  ' }
    43 = @{ D='0.779'; E='  -0.05%  ' }
    44 = @{ D='4.37'; E='  -1.28%  ' }
    45 = @{ E='  -2.48%  ' }
    46 = @{ E='  +1.11%  ' }
    47 = @{ D='2.445.84'; E='  -2.73%  ' }
    48 = @{ D='22.99'; E='  -0.04%  ' }
    49 = @{ D='6.63'; E='  -2.03%  ' }
    50 = @{ D='2.12'; E='  +9.73%  ' }
    51 = @{ E='  -2.26%  ' }
}

foreach ($row in $updates.Keys) {
    $cells = $updates[$row]
    if ($cells.ContainsKey("D")) {
        $addr = "D$row"
        $text = $cells["D"]
        $isNumberLike = $text -match "^[+-]?[0-9]*\.?[0-9]+$"
        if ($isNumberLike) {
            # Force text storage (prefix like typing '123 in Excel) so the
            # numeric-looking string is not auto-converted to a number, then
            # reset the cell style back to Normal so no stray quote-prefix
            # style remains attached to the cell (matches original formatting).
            $ws.Range($addr).Value = "'" + $text
            $ws.Range($addr).Style = 'Normal'
        } else {
            $ws.Range($addr).Value = $text
        }
    }
    if ($cells.ContainsKey("E")) {
        $addr = "E$row"
        $ws.Range($addr).Value = $cells["E"]
    }
}
